# "sending Projects module test cases"
#
# The workbook ships sample/test data for the Projects sheet: a header row
# (tyss / dinga / hello) and one data row (date, qty, amount). The Clients
# sheet - previously the active tab with an empty grid - loses focus to the
# Projects sheet, which becomes the new active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Clients")
$ws8 = $wb.Worksheets.Item("Projects")

# --- Clients sheet: no longer the active tab; gains a couple of touched,
#     but still content-free, rows and a new selection -------------------
$ws1.Range("A1:C2").Font.Bold = $false
$ws1.Range("B2").Select()

# --- Projects sheet: real test data -------------------------------------
# Column A is a bit wider to fit the date column.
$ws8.Columns.Item(1).ColumnWidth = 12.14

# Header row - written in this column order so the shared-string table
# ends up indexed tyss=0, hello=1, dinga=2.
$ws8.Range("A1").Value = "tyss"
$ws8.Range("C1").Value = "hello"
$ws8.Range("B1").Value = "dinga"

# Data row. Set the number format before the value so the engine doesn't
# auto-create a throwaway custom format and instead reuses the built-in
# "d-mmm-yy" (numFmtId 15).
$ws8.Range("A2").NumberFormat = "d-mmm-yy"
$ws8.Range("A2").Value = (Get-Date -Year 2021 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws8.Range("B2").Value = 2
$ws8.Range("C2").Value = 20.2

# Selection / active-sheet state.
$ws8.Range("B10").Select()
$ws8.Select()
